$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Info": update objective/time summary values
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 640108574274.0112
$wsInfo.Range("B2").Value = 1.670000076293945

# ---------------------------------------------------------------------------
# Sheet "Activados": Proceso/Tiempo values, extended from row 4 to row 20.
# Proceso column changes from 4 to 1, Tiempo now steps by 20 up to 360
# ---------------------------------------------------------------------------
$wsAct = $wb.Worksheets.Item("Activados")
for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $wsAct.Cells.Item($row, 1).Value = 1
    $wsAct.Cells.Item($row, 2).Value = $i * 20
}

# ---------------------------------------------------------------------------
# Sheet "Operando": Proceso column (A) changes from 4 to 1 for every data
# row (2-366); Tiempo column (B) stays the same
# ---------------------------------------------------------------------------
$wsOp = $wb.Worksheets.Item("Operando")
$wsOp.Range("A2:A366").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Contaminantes": update Z/Concentracion values for rows 2-6
# ---------------------------------------------------------------------------
$wsCont = $wb.Worksheets.Item("Contaminantes")

$wsCont.Cells.Item(2, 2).Value = 449208244800.0004
$wsCont.Cells.Item(2, 3).Value = 16.66000000000001

$wsCont.Cells.Item(3, 2).Value = 13481640000.00001
$wsCont.Cells.Item(3, 3).Value = 0.5000000000000004

$wsCont.Cells.Item(4, 2).Value = 87091394399.99998
$wsCont.Cells.Item(4, 3).Value = 3.23

$wsCont.Cells.Item(5, 2).Value = 307074.010608
$wsCont.Cells.Item(5, 3).Value = [double]"1.13886e-05"

$wsCont.Cells.Item(6, 2).Value = 90326988000.00008
$wsCont.Cells.Item(6, 3).Value = 3.350000000000003
